$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all cell contents (this also resets the shared-strings table so we can
# rebuild it below in the exact order the target workbook expects, while the
# ClearContents call preserves each cell's existing style index).
$ws.Cells.ClearContents()

# ---- Header row (row 1): add new "Label" column header in H1 ----
$ws.Cells.Item(1,2).Value = "Iterations"
$ws.Cells.Item(1,3).Value = "Success"
$ws.Cells.Item(1,4).Value = "Prediction"
$ws.Cells.Item(1,5).Value = "Error"
$ws.Cells.Item(1,6).Value = "Cross Entropy Loss"
$ws.Cells.Item(1,7).Value = "Success %"
$ws.Cells.Item(1,8).Value = "Label"

# H1 needs the same bold/centered/bordered header style as the rest of row 1.
# Copy formatting from an existing styled header cell rather than rebuilding
# the style piecemeal (keeps a single shared cellXf instead of duplicates).
$ws.Cells.Item(1,2).Copy()
$ws.Cells.Item(1,8).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Data rows ----
# Column A labels (re-introduced in row order so the shared-string table is
# rebuilt as: headers..., Label, Control 26, Control 33, Control 36,
# Control 49, Control 2, MDD 36, MDD 10, MDD 39, MDD 14, MDD 18)
$ws.Cells.Item(2,1).Value = "Control 26"
$ws.Cells.Item(3,1).Value = "Control 33"
$ws.Cells.Item(4,1).Value = "Control 36"
$ws.Cells.Item(5,1).Value = "Control 49"
$ws.Cells.Item(6,1).Value = "Control 2"
$ws.Cells.Item(7,1).Value = "MDD 36"
$ws.Cells.Item(8,1).Value = "MDD 10"
$ws.Cells.Item(9,1).Value = "MDD 39"
$ws.Cells.Item(10,1).Value = "MDD 14"
$ws.Cells.Item(11,1).Value = "MDD 18"
$ws.Cells.Item(12,1).Value = "Control 26"
$ws.Cells.Item(13,1).Value = "Control 33"
$ws.Cells.Item(14,1).Value = "Control 36"
$ws.Cells.Item(15,1).Value = "Control 49"
$ws.Cells.Item(16,1).Value = "Control 2"
$ws.Cells.Item(17,1).Value = "MDD 36"
$ws.Cells.Item(18,1).Value = "MDD 10"
$ws.Cells.Item(19,1).Value = "MDD 39"
$ws.Cells.Item(20,1).Value = "MDD 14"
$ws.Cells.Item(21,1).Value = "MDD 18"

# Column B (Iterations) - only set on the first row of each batch
$ws.Cells.Item(2,2).Value = 100
$ws.Cells.Item(12,2).Value = 200

# Column C (Success, boolean)
$ws.Cells.Item(2,3).Value = $true
$ws.Cells.Item(3,3).Value = $false
$ws.Cells.Item(4,3).Value = $false
$ws.Cells.Item(5,3).Value = $true
$ws.Cells.Item(6,3).Value = $false
$ws.Cells.Item(7,3).Value = $false
$ws.Cells.Item(8,3).Value = $true
$ws.Cells.Item(9,3).Value = $true
$ws.Cells.Item(10,3).Value = $true
$ws.Cells.Item(11,3).Value = $true
$ws.Cells.Item(12,3).Value = $true
$ws.Cells.Item(13,3).Value = $false
$ws.Cells.Item(14,3).Value = $false
$ws.Cells.Item(15,3).Value = $true
$ws.Cells.Item(16,3).Value = $false
$ws.Cells.Item(17,3).Value = $false
$ws.Cells.Item(18,3).Value = $true
$ws.Cells.Item(19,3).Value = $true
$ws.Cells.Item(20,3).Value = $true
$ws.Cells.Item(21,3).Value = $true

# Column D (Prediction) - refit values
$ws.Cells.Item(2,4).Value = 0.4742328669449746
$ws.Cells.Item(3,4).Value = 0.6691163950210938
$ws.Cells.Item(4,4).Value = 0.5657883523574412
$ws.Cells.Item(5,4).Value = 0.3479351807662002
$ws.Cells.Item(6,4).Value = 0.6359048007071821
$ws.Cells.Item(7,4).Value = 0.4741486656035482
$ws.Cells.Item(8,4).Value = 0.6119415999995075
$ws.Cells.Item(9,4).Value = 0.8036737623713639
$ws.Cells.Item(10,4).Value = 0.7010507576668658
$ws.Cells.Item(11,4).Value = 0.7655308865889124
$ws.Cells.Item(12,4).Value = 0.4813196907902212
$ws.Cells.Item(13,4).Value = 0.6691163950210938
$ws.Cells.Item(14,4).Value = 0.5653844334601101
$ws.Cells.Item(15,4).Value = 0.3479351807662002
$ws.Cells.Item(16,4).Value = 0.6294168126217071
$ws.Cells.Item(17,4).Value = 0.4813113165951227
$ws.Cells.Item(18,4).Value = 0.6075592283222935
$ws.Cells.Item(19,4).Value = 0.8017805534951261
$ws.Cells.Item(20,4).Value = 0.6947236202701031
$ws.Cells.Item(21,4).Value = 0.7620974157734631

# Column E (Error) - refit values
$ws.Cells.Item(2,5).Value = 0.4742328669449746
$ws.Cells.Item(3,5).Value = 0.6691163950210938
$ws.Cells.Item(4,5).Value = 0.5657883523574412
$ws.Cells.Item(5,5).Value = 0.3479351807662002
$ws.Cells.Item(6,5).Value = 0.6359048007071821
$ws.Cells.Item(7,5).Value = 0.5258513343964517
$ws.Cells.Item(8,5).Value = 0.3880584000004925
$ws.Cells.Item(9,5).Value = 0.1963262376286361
$ws.Cells.Item(10,5).Value = 0.2989492423331342
$ws.Cells.Item(11,5).Value = 0.2344691134110876
$ws.Cells.Item(12,5).Value = 0.4813196907902212
$ws.Cells.Item(13,5).Value = 0.6691163950210938
$ws.Cells.Item(14,5).Value = 0.5653844334601101
$ws.Cells.Item(15,5).Value = 0.3479351807662002
$ws.Cells.Item(16,5).Value = 0.6294168126217071
$ws.Cells.Item(17,5).Value = 0.5186886834048773
$ws.Cells.Item(18,5).Value = 0.3924407716777065
$ws.Cells.Item(19,5).Value = 0.1982194465048739
$ws.Cells.Item(20,5).Value = 0.3052763797298969
$ws.Cells.Item(21,5).Value = 0.2379025842265369

# Columns F & G (Cross Entropy Loss / Success %) - only present for the last
# row of each batch
$ws.Cells.Item(11,6).Value = 0.6099335551261902
$ws.Cells.Item(11,7).Value = 0.6
$ws.Cells.Item(21,6).Value = 0.6102427840232849
$ws.Cells.Item(21,7).Value = 0.6

# Column H (new "Label" column: 0 = Control, 1 = MDD)
$ws.Cells.Item(2,8).Value = 0
$ws.Cells.Item(3,8).Value = 0
$ws.Cells.Item(4,8).Value = 0
$ws.Cells.Item(5,8).Value = 0
$ws.Cells.Item(6,8).Value = 0
$ws.Cells.Item(7,8).Value = 1
$ws.Cells.Item(8,8).Value = 1
$ws.Cells.Item(9,8).Value = 1
$ws.Cells.Item(10,8).Value = 1
$ws.Cells.Item(11,8).Value = 1
$ws.Cells.Item(12,8).Value = 0
$ws.Cells.Item(13,8).Value = 0
$ws.Cells.Item(14,8).Value = 0
$ws.Cells.Item(15,8).Value = 0
$ws.Cells.Item(16,8).Value = 0
$ws.Cells.Item(17,8).Value = 1
$ws.Cells.Item(18,8).Value = 1
$ws.Cells.Item(19,8).Value = 1
$ws.Cells.Item(20,8).Value = 1
$ws.Cells.Item(21,8).Value = 1
